$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 22; $row++) {
    if ($row -eq 2) {
        $ws.Cells.Item($row, 2).Value = 251442
    } else {
        $ws.Cells.Item($row, 2).Value = 248232
    }

    if ($row -eq 4) {
        $ws.Cells.Item($row, 3).Value = 271234
    } else {
        $ws.Cells.Item($row, 3).Value = 275600
    }

    $ws.Cells.Item($row, 4).Value = 272396
}
